$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.889.25"
Set-TextValue $ws.Range("D3") "1.545.07"
Set-TextValue $ws.Range("E3") "  -1.20%  "
Set-TextValue $ws.Range("E4") "  +0.22%  "
Set-TextValue $ws.Range("D5") "205.81"
Set-TextValue $ws.Range("E5") "  -0.06%  "
Set-TextValue $ws.Range("D6") "0.484"
Set-TextValue $ws.Range("E6") "  -0.68%  "
Set-TextValue $ws.Range("E7") "  +0.26%  "
Set-TextValue $ws.Range("D9") "21.27"
Set-TextValue $ws.Range("E9") "  -2.28%  "
Set-TextValue $ws.Range("E10") "  -0.43%  "
Set-TextValue $ws.Range("D11") "0.0856"
Set-TextValue $ws.Range("E11") "  -0.95%  "
Set-TextValue $ws.Range("D12") "1.764.88"
Set-TextValue $ws.Range("E12") "  -1.19%  "
Set-TextValue $ws.Range("D13") "1.540.80"
Set-TextValue $ws.Range("E13") "  -1.42%  "
Set-TextValue $ws.Range("E14") "  -0.93%  "
Set-TextValue $ws.Range("D15") "0.510"
Set-TextValue $ws.Range("E15") "  -0.93%  "
Set-TextValue $ws.Range("D16") "26.863.52"
Set-TextValue $ws.Range("D17") "61.46"
Set-TextValue $ws.Range("E17") "  +0.33%  "
Set-TextValue $ws.Range("D18") "213.24"
Set-TextValue $ws.Range("E18") "  -0.80%  "
Set-TextValue $ws.Range("E19") "  +0.12%  "
Set-TextValue $ws.Range("D20") "7.19"
Set-TextValue $ws.Range("E20") "  -2.32%  "
Set-TextValue $ws.Range("E21") "  +0.35%  "
Set-TextValue $ws.Range("E22") "  -2.58%  "
Set-TextValue $ws.Range("D23") "9.19"
Set-TextValue $ws.Range("E23") "  -0.15%  "
Set-TextValue $ws.Range("D25") "152.81"
Set-TextValue $ws.Range("E25") "  -0.42%  "
Set-TextValue $ws.Range("E26") "  -1.77%  "
Set-TextValue $ws.Range("D27") "14.81"
Set-TextValue $ws.Range("E27") "  -0.69%  "
Set-TextValue $ws.Range("E28") "  +0.27%  "
Set-TextValue $ws.Range("E29") "  -0.13%  "
Set-TextValue $ws.Range("E30") "  -1.88%  "
Set-TextValue $ws.Range("E31") "  -0.86%  "
Set-TextValue $ws.Range("E32") "  +1.43%  "
Set-TextValue $ws.Range("D33") "1.357.62"
Set-TextValue $ws.Range("E33") "  -3.43%  "
Set-TextValue $ws.Range("E34") "  +0.45%  "
Set-TextValue $ws.Range("D35") "1.53"
Set-TextValue $ws.Range("E35") "  +0.09%  "
Set-TextValue $ws.Range("E36") "  +5.18%  "
Set-TextValue $ws.Range("E37") "  +0.33%  "
Set-TextValue $ws.Range("E38") "  +0.09%  "
Set-TextValue $ws.Range("D39") "0.519"
Set-TextValue $ws.Range("E39") "  -1.58%  "
Set-TextValue $ws.Range("E40") "  -0.74%  "
Set-TextValue $ws.Range("E41") "  +0.23%  "
Set-TextValue $ws.Range("D42") "5.61"
Set-TextValue $ws.Range("E42") "  +2.76%  "
Set-TextValue $ws.Range("E43") "  -0.72%  "
Set-TextValue $ws.Range("D44") "2.21"
Set-TextValue $ws.Range("E44") "  +1.75%  "
Set-TextValue $ws.Range("D45") "63.46"
Set-TextValue $ws.Range("E45") "  +0.19%  "
Set-TextValue $ws.Range("E46") "  -2.27%  "
Set-TextValue $ws.Range("D47") "1.678.94"
Set-TextValue $ws.Range("E47") "  -1.28%  "
Set-TextValue $ws.Range("D48") "85.96"
Set-TextValue $ws.Range("E48") "  -0.73%  "
Set-TextValue $ws.Range("E49") "  +0.79%  "
Set-TextValue $ws.Range("E50") "  -0.90%  "
Set-TextValue $ws.Range("D51") "0.0946"
Set-TextValue $ws.Range("E51") "  -0.14%  "
